$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# values such as "1.001" or "30.145.52" are not re-interpreted as numbers.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D25","D26","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, row by row, in the same order as the diff.
$ws.Range("D2").Value = '30.145.52'
$ws.Range("E2").Value = '  +5.61%  '
$ws.Range("D3").Value = '1.927.32'
$ws.Range("E3").Value = '  +3.03%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.70%  '
$ws.Range("D5").Value = '322.95'
$ws.Range("E5").Value = '  +2.64%  '
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("D7").Value = '0.5162'
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("D8").Value = '0.3991'
$ws.Range("E8").Value = '  +2.47%  '
$ws.Range("D9").Value = '0.08489'
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("D10").Value = '43.02'
$ws.Range("E10").Value = '  +3.02%  '
$ws.Range("E11").Value = '  +2.13%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '6.334'
$ws.Range("E12").Value = '  +2.06%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '21.18'
$ws.Range("E13").Value = '  +4.02%  '
$ws.Range("D14").Value = '1.923.92'
$ws.Range("E14").Value = '  +2.78%  '
$ws.Range("D15").Value = '7.384'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("D17").Value = '94.44'
$ws.Range("E17").Value = '  +3.84%  '
$ws.Range("D18").Value = '0.00001119'
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("D19").Value = '0.06757'
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("D20").Value = '18.01'
$ws.Range("E20").Value = '  +2.03%  '
$ws.Range("D21").Value = '0.9996'
$ws.Range("E21").Value = '  -0.79%  '
$ws.Range("D22").Value = '6.080'
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").Value = '30.185.21'
$ws.Range("E23").Value = '  +5.64%  '
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("D25").Value = '2.208'
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("D26").Value = '2.140.49'
$ws.Range("E26").Value = '  +2.60%  '
$ws.Range("E27").Value = '  +2.03%  '
$ws.Range("D28").Value = '159.60'
$ws.Range("E28").Value = '  -1.29%  '
$ws.Range("D29").Value = '2.476'
$ws.Range("E29").Value = '  +6.00%  '
$ws.Range("D30").Value = '129.11'
$ws.Range("E30").Value = '  +2.70%  '
$ws.Range("D31").Value = '1.081'
$ws.Range("E31").Value = '  +3.93%  '
$ws.Range("D32").Value = '0.1059'
$ws.Range("E32").Value = '  +1.69%  '
$ws.Range("D33").Value = '6.106'
$ws.Range("E33").Value = '  +5.67%  '
$ws.Range("D34").Value = '3.658'
$ws.Range("E34").Value = '  +1.27%  '
$ws.Range("D35").Value = '0.02501'
$ws.Range("E35").Value = '  +2.29%  '
$ws.Range("D36").Value = '0.06634'
$ws.Range("E36").Value = '  +1.68%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2213'
$ws.Range("E37").Value = '  +2.52%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = '1.246'
$ws.Range("E38").Value = '  +5.04%  '
$ws.Range("D39").Value = '9.075'
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("D40").Value = '5.202'
$ws.Range("E40").Value = '  +3.17%  '
$ws.Range("D41").Value = '0.6542'
$ws.Range("E41").Value = '  +1.94%  '
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("D43").Value = '11.41'
$ws.Range("E43").Value = '  +2.93%  '
$ws.Range("D44").Value = '0.6153'
$ws.Range("E44").Value = '  +2.32%  '
$ws.Range("D45").Value = '13.24'
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("D46").Value = '3.727'
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("D47").Value = '2.060'
$ws.Range("E47").Value = '  +2.98%  '
$ws.Range("D48").Value = '1.243'
$ws.Range("E48").Value = '  +2.45%  '
$ws.Range("D49").Value = '125.56'
$ws.Range("E49").Value = '  +3.11%  '
$ws.Range("D50").Value = '1.150'
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("D51").Value = '79.43'
$ws.Range("E51").Value = '  +4.17%  '
